# Updated cryptos list on Sun Dec  1 17:41:56 UTC 2024 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns on Sheet1 with the latest
# scrape, including two rank swaps (rows 29/30 and 36/37) where the underlying
# coins traded places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "97.218.85"
$ws.Range("E2").Value = "  +0.35%  "

# Row 3
$ws.Range("D3").Value = "3.712.05"
$ws.Range("E3").Value = "  +0.84%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.87"
$ws.Range("E5").Value = "  -1.15%  "

# Row 6
$ws.Range("E6").Value = "  +1.48%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "657.89"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.433"
$ws.Range("E8").Value = "  +1.83%  "

# Row 9
$ws.Range("E9").Value = "  -0.03%  "

# Row 10
$ws.Range("E10").Value = "  -2.21%  "

# Row 11
$ws.Range("D11").Value = "3.709.25"
$ws.Range("E11").Value = "  +0.83%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000318"
$ws.Range("E12").Value = "  +17.43%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "44.66"
$ws.Range("E13").Value = "  -1.96%  "

# Row 14
$ws.Range("E14").Value = "  +0.83%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.95"
$ws.Range("E15").Value = "  +1.95%  "

# Row 16
$ws.Range("D16").Value = "4.405.73"
$ws.Range("E16").Value = "  +0.90%  "

# Row 17
$ws.Range("D17").Value = "96.903.56"
$ws.Range("E17").Value = "  +0.36%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.00"
$ws.Range("E18").Value = "  +0.85%  "

# Row 19
$ws.Range("D19").Value = "3.702.88"
$ws.Range("E19").Value = "  +0.71%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.00"
$ws.Range("E20").Value = "  +1.97%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.69"
$ws.Range("E21").Value = "  -0.52%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.508"
$ws.Range("E22").Value = "  -4.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "525.48"
$ws.Range("E23").Value = "  -1.09%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.48"
$ws.Range("E24").Value = "  -1.00%  "

# Row 25
$ws.Range("E25").Value = "  +9.24%  "

# Row 26
$ws.Range("E26").Value = "  -3.58%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "105.89"
$ws.Range("E27").Value = "  +3.16%  "

# Row 28
$ws.Range("E28").Value = "  +0.93%  "

# Row 29
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.55"
$ws.Range("E29").Value = "  +0.24%  "

# Row 30
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.189"
$ws.Range("E30").Value = "  +12.92%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.67"
$ws.Range("E31").Value = "  +2.45%  "

# Row 32
$ws.Range("E32").Value = "  -0.93%  "

# Row 34
$ws.Range("E34").Value = "  +3.15%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.84"
$ws.Range("E35").Value = "  -4.15%  "

# Row 36
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.51"
$ws.Range("E36").Value = "  -0.49%  "

# Row 37
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.23%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "639.17"
$ws.Range("E38").Value = "  -3.62%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.591"
$ws.Range("E39").Value = "  -1.08%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.77"
$ws.Range("E40").Value = "  -1.34%  "

# Row 41
$ws.Range("E41").Value = "  +0.03%  "

# Row 42
$ws.Range("E42").Value = "  +2.86%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.61"
$ws.Range("E43").Value = "  +4.78%  "

# Row 44
$ws.Range("E44").Value = "  +1.70%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.02"
$ws.Range("E45").Value = "  +0.94%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.481"
$ws.Range("E46").Value = "  +11.99%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.968"
$ws.Range("E47").Value = "  +0.37%  "

# Row 48
$ws.Range("E48").Value = "  -1.25%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.38"
$ws.Range("E49").Value = "  +1.95%  "

# Row 50
$ws.Range("E50").Value = "  -0.11%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.66"
$ws.Range("E51").Value = "  -0.30%  "
